# Update "想去人数" (want-to-go count) figures in column F across the four
# sheets of the workbook to reflect the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "展览" (Exhibition)
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 63
$ws.Range("F5").Value  = 79
$ws.Range("F6").Value  = 859
$ws.Range("F7").Value  = 439
$ws.Range("F8").Value  = 4751
$ws.Range("F9").Value  = 4751
$ws.Range("F11").Value = 125
$ws.Range("F14").Value = 200
$ws.Range("F15").Value = 125
$ws.Range("F16").Value = 7597
$ws.Range("F17").Value = 251
$ws.Range("F21").Value = 533
$ws.Range("F22").Value = 1397
$ws.Range("F24").Value = 6287
$ws.Range("F25").Value = 2257
$ws.Range("F28").Value = 6196
$ws.Range("F29").Value = 146
$ws.Range("F30").Value = 25
$ws.Range("F31").Value = 118
$ws.Range("F33").Value = 449
$ws.Range("F34").Value = 6489
$ws.Range("F35").Value = 25
$ws.Range("F36").Value = 211
$ws.Range("F40").Value = 32
$ws.Range("F41").Value = 2460
$ws.Range("F45").Value = 40
$ws.Range("F46").Value = 450
$ws.Range("F47").Value = 2152
$ws.Range("F49").Value = 1083

# -----------------------------------------------------------------
# Sheet "演出" (Performance)
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value  = 233
$ws.Range("F6").Value  = 129
$ws.Range("F13").Value = 148

# -----------------------------------------------------------------
# Sheet "本地生活" (Local life)
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1451

# -----------------------------------------------------------------
# Sheet "全部类型" (All types)
# -----------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 1451
$ws.Range("F4").Value  = 63
$ws.Range("F5").Value  = 233
$ws.Range("F6").Value  = 79
$ws.Range("F8").Value  = 439
$ws.Range("F9").Value  = 4751
$ws.Range("F10").Value = 4751
$ws.Range("F12").Value = 125
$ws.Range("F15").Value = 200
$ws.Range("F16").Value = 125
$ws.Range("F17").Value = 7597
$ws.Range("F18").Value = 251
$ws.Range("F20").Value = 533
$ws.Range("F21").Value = 1397
$ws.Range("F22").Value = 129
$ws.Range("F23").Value = 6287
$ws.Range("F24").Value = 2257
$ws.Range("F29").Value = 6196
$ws.Range("F30").Value = 146
$ws.Range("F32").Value = 25
$ws.Range("F33").Value = 118
$ws.Range("F35").Value = 449
$ws.Range("F36").Value = 6489
$ws.Range("F37").Value = 25
$ws.Range("F38").Value = 211
$ws.Range("F42").Value = 2460
$ws.Range("F45").Value = 40
$ws.Range("F46").Value = 450
$ws.Range("F47").Value = 148
$ws.Range("F48").Value = 2152

$wb.Save()
